$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (Price / Volume) to Text format first so that
# numeric-looking strings (e.g. "19.53", "1.00") are not silently
# reinterpreted by Excel as numbers, matching the source inlineStr cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.210.10'
$ws.Range("E2").Value = '  -1.98%  '
$ws.Range("D3").Value = '1.582.97'
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("E6").Value = '  -2.83%  '
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.246'
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.0610'
$ws.Range("E9").Value = '  -1.46%  '
$ws.Range("D10").Value = '19.53'
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("D11").Value = '0.0846'
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("D12").Value = '1.805.70'
$ws.Range("E12").Value = '  -1.15%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.595.48'
$ws.Range("E13").Value = '  -0.40%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '4.06'
$ws.Range("E14").Value = '  +0.49%  '
$ws.Range("D15").Value = '0.515'
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").Value = '64.52'
$ws.Range("E16").Value = '  -0.78%  '
$ws.Range("D17").Value = '26.215.46'
$ws.Range("E17").Value = '  -1.83%  '
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").Value = '7.29'
$ws.Range("E19").Value = '  +1.40%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '1.00'
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '207.31'
$ws.Range("E21").Value = '  -1.46%  '
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("D23").Value = '2.20'
$ws.Range("E23").Value = '  -3.08%  '
$ws.Range("E24").Value = '  -0.96%  '
$ws.Range("D25").Value = '144.47'
$ws.Range("E25").Value = '  +0.53%  '
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("D27").Value = '7.01'
$ws.Range("E27").Value = '  -1.26%  '
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("D29").Value = '15.23'
$ws.Range("E29").Value = '  -1.04%  '
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("D33").Value = '2.95'
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").Value = '1.277.03'
$ws.Range("E34").Value = '  -1.30%  '
$ws.Range("D35").Value = '2.47'
$ws.Range("E35").Value = '  -0.38%  '
$ws.Range("D36").Value = '0.613'
$ws.Range("E36").Value = '  +1.43%  '
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("D38").Value = '0.0166'
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").Value = '1.05'
$ws.Range("E39").Value = '  -9.98%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '0.818'
$ws.Range("E40").Value = '  -1.53%  '
$ws.Range("D41").Value = '5.55'
$ws.Range("E41").Value = '  +2.60%  '
$ws.Range("D42").Value = '0.767'
$ws.Range("E42").Value = '  -2.25%  '
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("D44").Value = '62.34'
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("D45").Value = '1.719.22'
$ws.Range("E45").Value = '  -1.17%  '
$ws.Range("D46").Value = '89.19'
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").Value = '0.100'
$ws.Range("E48").Value = '  -0.72%  '
$ws.Range("D49").Value = '0.0507'
$ws.Range("E49").Value = '  -1.98%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.45'
$ws.Range("E51").Value = '  +0.80%  '
